# SIPOT report workbook update: roll the reporting period forward from
# Q2 2021 to Q4 2021 ("4to Trimestre"), per the commit message
# "actualizacion de febrero hay un archivo mal".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reporte de Formatos")

# --- Page setup: paper size ---------------------------------------------
$ws.PageSetup.PaperSize = 9

# --- Data validations: shrink sqref from row 201 down to row 90 --------
$validations = @(
    @{ Col = "F"; Name = "Hidden_15" },
    @{ Col = "J"; Name = "Hidden_29" },
    @{ Col = "Q"; Name = "Hidden_316" },
    @{ Col = "W"; Name = "Hidden_422" },
    @{ Col = "X"; Name = "Hidden_523" },
    @{ Col = "Y"; Name = "Hidden_624" }
)
foreach ($v in $validations) {
    $oldRange = $ws.Range($v.Col + "8:" + $v.Col + "201")
    $oldRange.Validation.Delete()
    $newRange = $ws.Range($v.Col + "8:" + $v.Col + "90")
    $newRange.Validation.Add(3, 1, 1, $v.Name)
    $newRange.Validation.ShowInput = $false
}

# --- Row 8: new reporting period (Q4 2021) ------------------------------
$ws.Range("B8").Value2 = 44378   # 2021-07-01
$ws.Range("C8").Value2 = 44561   # 2021-12-31
$ws.Range("AG8").Value2 = 44571  # 2022-01-10
$ws.Range("AH8").Value2 = 44571  # 2022-01-10

# --- Selection -----------------------------------------------------------
$ws.Range("E10").Select()
